# Scraper ran again and picked up more data: column F (count of entries
# actually scraped) now has values for rows 27-41, which previously sat
# empty/pending (red "Bad" style). Filling them in flips their style to
# the green "Good" style, same as every already-completed row above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> newly scraped count
$newCounts = [ordered]@{
    27 = 14
    28 = 16
    29 = 1
    30 = 27
    31 = 2
    32 = 9
    33 = 14
    34 = 22
    35 = 12
    36 = 1
    37 = 1
    38 = 8
    39 = 9
    40 = 18
    41 = 14
}

# F2 already carries the "complete" (Good) look these rows need once they
# have data. Copy its formatting onto each newly-completed row first ...
$ws.Range("F2").Copy()
foreach ($row in $newCounts.Keys) {
    $ws.Range("F$row").PasteSpecial(-4122)
}

# ... then stamp in the actual scraped counts.
foreach ($row in $newCounts.Keys) {
    $ws.Range("F$row").Value = $newCounts[$row]
}

# Leave the selection where the author's cursor ended up after scrolling
# through the newly-populated rows.
$ws.Range("F42").Select()
